# Estadisticos Segundo Parcial 23 Mayo
# This workbook tracks pass/fail statistics per subject (Mat) across
# three grading periods (1P, 2P, Final) and a "Rescatables" (make-up)
# roster. The edit reflects a re-grade: one student who previously failed
# "IMPLEMENTA BASE DE DATOS..." now passes, which shifts the aggregate
# counts/percentages/averages on all three stats sheets and updates the
# make-up roster to show the new student in need of a retake there
# (ORTIZ CORTES CARLOS) while removing the now-passed student for that
# subject and dropping all the other stale make-up rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Estadisticos 1P"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Cells.Item(3,3).Value = 0      # C3 Totales
$ws1.Cells.Item(3,6).Value = 0      # F3 Aprobados
$ws1.Cells.Item(3,7).ClearContents()  # G3 Por_Apro removed
$ws1.Cells.Item(3,8).ClearContents()  # H3 Promedio removed

# ---------------------------------------------------------------
# Sheet "Estadisticos 2P"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Cells.Item(2,4).Value = 0        # D2 Blancos
$ws2.Cells.Item(2,5).Value = 4        # E2 Reprobados
$ws2.Cells.Item(2,6).Value = 19       # F2 Aprobados
$ws2.Cells.Item(2,7).Value = 82.61    # G2 Por_Apro
$ws2.Cells.Item(2,8).Value = 7.6      # H2 Promedio (new)

$ws2.Cells.Item(3,3).Value = 0        # C3 Totales
$ws2.Cells.Item(3,4).Value = 0        # D3 Blancos
$ws2.Cells.Item(3,5).Value = 0        # E3 Reprobados
$ws2.Cells.Item(3,7).ClearContents()  # G3 Por_Apro removed

$ws2.Cells.Item(4,4).Value = 0        # D4 Blancos
$ws2.Cells.Item(4,5).Value = 2        # E4 Reprobados
$ws2.Cells.Item(4,6).Value = 11       # F4 Aprobados
$ws2.Cells.Item(4,7).Value = 84.62    # G4 Por_Apro
$ws2.Cells.Item(4,8).Value = 8.1      # H4 Promedio (new)

# ---------------------------------------------------------------
# Sheet "Estadisticos Final"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Cells.Item(2,5).Value = 4          # E2 Reprobados
$ws3.Cells.Item(2,6).Value = 19         # F2 Aprobados
$ws3.Cells.Item(2,7).Value = 82.61      # G2 Por_Apro
$ws3.Cells.Item(2,8).Value = 7.8        # H2 Promedio

$ws3.Cells.Item(3,3).Value = 0          # C3 Totales
$ws3.Cells.Item(3,4).Value = 0          # D3 Blancos
$ws3.Cells.Item(3,7).ClearContents()    # G3 Por_Apro removed

$ws3.Cells.Item(4,8).Value = 8.199999999999999   # H4 Promedio

# ---------------------------------------------------------------
# Sheet "Rescatables"
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

# Drop the stale rows 4-13; only the header + two roster rows remain.
$ws4.Range("A4:G13").EntireRow.Delete()

# Row 2: new "rescatable" student for IMPLEMENTA BASE DE DATOS...
$ws4.Cells.Item(2,1).Value = 23330051920301
$ws4.Cells.Item(2,2).Value = "ORTIZ"
$ws4.Cells.Item(2,3).Value = "CORTES"
$ws4.Cells.Item(2,4).Value = "CARLOS"
$ws4.Cells.Item(2,5).Value = "IMPLEMENTA BASE DE DATOS RELACIONALES EN UN SISTEMA DE INFORMACIÓN"
$ws4.Cells.Item(2,6).Value = "4APV"
$ws4.Cells.Item(2,7).Value = 4

# Row 3: previously-failing student, now also IMPLEMENTA BASE DE DATOS...
$ws4.Cells.Item(3,1).Value = 23330051920317
$ws4.Cells.Item(3,2).Value = "ESTRADA"
$ws4.Cells.Item(3,3).Value = "SAN JUAN"
$ws4.Cells.Item(3,4).Value = "JOSE MARCOS"
$ws4.Cells.Item(3,5).Value = "IMPLEMENTA BASE DE DATOS RELACIONALES EN UN SISTEMA DE INFORMACIÓN"
$ws4.Cells.Item(3,6).Value = "4APV"
$ws4.Cells.Item(3,7).Value = 3
